# ---- Rename sheets: "zz" -> "s1", "xx" -> "s2" ----
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Name = "s1"
$ws2.Name = "s2"

# ---- Populate the (previously empty) second sheet, "s2", with a
#      service-item before/after comparison table (34 cols x 3 rows) ----
$headerVals = New-Object 'object[,]' 1,34
$headerVals[0,0] = "Service Item ID"
$headerVals[0,1] = "Item Type"
$headerVals[0,2] = "Folder"
$headerVals[0,3] = "Folder (after)"
$headerVals[0,4] = "Item Title"
$headerVals[0,5] = "Item Title (after)"
$headerVals[0,6] = "Shared With"
$headerVals[0,7] = "Shared With (after)"
$headerVals[0,8] = "Layer Title"
$headerVals[0,9] = "Layer Title (after)"
$headerVals[0,10] = "Layer Visibility"
$headerVals[0,11] = "Layer Visibility (after)"
$headerVals[0,12] = "Layer Item ID"
$headerVals[0,13] = "Layer Item ID (after)"
$headerVals[0,14] = "Popup Title"
$headerVals[0,15] = "Popup Title (after)"
$headerVals[0,16] = "Form Configured"
$headerVals[0,17] = "Form Configured (after)"
$headerVals[0,18] = "FormField"
$headerVals[0,19] = "FormField (after)"
$headerVals[0,20] = "FormLabel"
$headerVals[0,21] = "FormLabel (after)"
$headerVals[0,22] = "FormFieldEditable"
$headerVals[0,23] = "FormFieldEditable (after)"
$headerVals[0,24] = "FormFieldInputType"
$headerVals[0,25] = "FormFieldInputType (after)"
$headerVals[0,26] = "Calculated Expression"
$headerVals[0,27] = "Calculated Expression (after)"
$headerVals[0,28] = "FormFieldRequired"
$headerVals[0,29] = "FormFieldRequired (after)"
$headerVals[0,30] = "Offline Status"
$headerVals[0,31] = "Offline Status (after)"
$headerVals[0,32] = "Delete Protection"
$headerVals[0,33] = "Delete Protection (after)"

$dataVals = New-Object 'object[,]' 2,34
$dataVals[0,0] = "c7654927cf554c9490fe639178351fb9"
$dataVals[0,1] = "Web Map"
$dataVals[0,2] = "ARM 2023 Space"
$dataVals[0,3] = "ARM 2023 Space Test"
$dataVals[0,4] = "LP_TestMap"
$dataVals[0,5] = "LP_TestMap"
$dataVals[0,6] = $null
$dataVals[0,7] = $null
$dataVals[0,8] = "LandPlanner_Area_Dev2023"
$dataVals[0,9] = "LandPlanner_Area_Dev2023"
$dataVals[0,10] = $true
$dataVals[0,11] = $true
$dataVals[0,12] = "c49a708cc0ac4ca1a634228ca7de0bf6"
$dataVals[0,13] = "c49a708cc0ac4ca1a634228ca7de0bf6"
$dataVals[0,14] = "{OBJECTID}"
$dataVals[0,15] = "{OBJECTID}"
$dataVals[0,16] = "Not Configured"
$dataVals[0,17] = "Not Configured"
$dataVals[0,18] = $null
$dataVals[0,19] = $null
$dataVals[0,20] = $null
$dataVals[0,21] = $null
$dataVals[0,22] = $null
$dataVals[0,23] = $null
$dataVals[0,24] = $null
$dataVals[0,25] = $null
$dataVals[0,26] = $null
$dataVals[0,27] = $null
$dataVals[0,28] = $null
$dataVals[0,29] = $null
$dataVals[0,30] = $null
$dataVals[0,31] = $null
$dataVals[0,32] = "On"
$dataVals[0,33] = "Off"
$dataVals[1,0] = "6003d30a5f06445897ac7facd7287480"
$dataVals[1,1] = "Web Map"
$dataVals[1,2] = "ARM_b Transmission Roads"
$dataVals[1,3] = "ARM_b Transmission Roads"
$dataVals[1,4] = "All Electric Transmission Operations Map"
$dataVals[1,5] = "All Electric Transmission Operations Map"
$dataVals[1,6] = "Test"
$dataVals[1,7] = $null
$dataVals[1,8] = "Workcenter"
$dataVals[1,9] = "Workcenter"
$dataVals[1,10] = $false
$dataVals[1,11] = $false
$dataVals[1,12] = "c1c831c25b44444b9491df15c56ce7bc"
$dataVals[1,13] = "c1c831c25b44444b9491df15c56ce7bc"
$dataVals[1,14] = "Workcenter: {WORKCENTER_ARM}"
$dataVals[1,15] = "Workcenter: {WORKCENTER_ARM}"
$dataVals[1,16] = "Not Configured"
$dataVals[1,17] = "Not Configured"
$dataVals[1,18] = $null
$dataVals[1,19] = $null
$dataVals[1,20] = $null
$dataVals[1,21] = $null
$dataVals[1,22] = $null
$dataVals[1,23] = $null
$dataVals[1,24] = $null
$dataVals[1,25] = $null
$dataVals[1,26] = $null
$dataVals[1,27] = $null
$dataVals[1,28] = $null
$dataVals[1,29] = $null
$dataVals[1,30] = $null
$dataVals[1,31] = $null
$dataVals[1,32] = "Off"
$dataVals[1,33] = "Off"

$ws2.Range("A1:AH1").Value = $headerVals
$ws2.Range("A2:AH3").Value = $dataVals

# ---- Header row formatting: bold, centered, thin box border (mirrors "s1") ----
$headerRange = $ws2.Range("A1:AH1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# ---- Highlight the before/after pairs whose values actually changed ----
$ws2.Range("C2:D2").Interior.Color = 65535
$ws2.Range("AG2:AH2").Interior.Color = 65535
$ws2.Range("G3:H3").Interior.Color = 65535

Write-Output "Sheets renamed and comparison data written to s2."
